# Update "想去人数" (interest count) and a couple "最低票价" (min price) values
# across the four worksheets, matching the upstream data refresh.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 (Exhibitions) ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F3").Value  = 361
$ws1.Range("F4").Value  = 413
$ws1.Range("F5").Value  = 1122
$ws1.Range("F8").Value  = 895
$ws1.Range("F9").Value  = 1609
$ws1.Range("F10").Value = 6064
$ws1.Range("G10").Value = 68
$ws1.Range("F11").Value = 111
$ws1.Range("F12").Value = 1740
$ws1.Range("F13").Value = 442
$ws1.Range("F14").Value = 5905
$ws1.Range("F15").Value = 114
$ws1.Range("F19").Value = 1646
$ws1.Range("F20").Value = 837
$ws1.Range("F22").Value = 140
$ws1.Range("F23").Value = 1369
$ws1.Range("F24").Value = 721
$ws1.Range("F25").Value = 237
$ws1.Range("F30").Value = 3858

# --- Sheet 2: 演出 (Performances) ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F4").Value  = 309
$ws2.Range("F5").Value  = 159
$ws2.Range("F20").Value = 28

# --- Sheet 3: 本地生活 (Local life) ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("F3").Value = 2228
$ws3.Range("F4").Value = 615
$ws3.Range("F5").Value = 178

# --- Sheet 4: 全部类型 (All types, combined listing) ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F3").Value  = 2228
$ws4.Range("F4").Value  = 615
$ws4.Range("F5").Value  = 361
$ws4.Range("F6").Value  = 413
$ws4.Range("F7").Value  = 1122
$ws4.Range("F11").Value = 309
$ws4.Range("F12").Value = 895
$ws4.Range("F13").Value = 178
$ws4.Range("F14").Value = 1609
$ws4.Range("F15").Value = 6064
$ws4.Range("G15").Value = 68
$ws4.Range("F16").Value = 111
$ws4.Range("F17").Value = 1740
$ws4.Range("F20").Value = 442
$ws4.Range("F23").Value = 5905
$ws4.Range("F24").Value = 114
$ws4.Range("F28").Value = 1646
$ws4.Range("F29").Value = 837
$ws4.Range("F31").Value = 140
$ws4.Range("F32").Value = 1369
$ws4.Range("F33").Value = 721
$ws4.Range("F35").Value = 237
$ws4.Range("F45").Value = 3858
$ws4.Range("F46").Value = 28
